# Updates cryptos list figures (price + 1h volume change) to match the latest
# scrape, and fixes the ApeXProtocol / ThetaToken row ordering (rows 46-47 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.565.66"
$ws.Range("E2").Value = "'  +3.27%  "
$ws.Range("D3").Value = "'3.709.65"
$ws.Range("E3").Value = "'  +8.60%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'582.87"
$ws.Range("E5").Value = "'  +0.44%  "
$ws.Range("D6").Value = "'178.75"
$ws.Range("E6").Value = "'  +1.42%  "
$ws.Range("D7").Value = "'3.699.18"
$ws.Range("E7").Value = "'  +8.48%  "
$ws.Range("E8").Value = "'  +4.39%  "
$ws.Range("D10").Value = "'0.200"
$ws.Range("E10").Value = "'  +1.96%  "
$ws.Range("E11").Value = "'  +4.91%  "
$ws.Range("D12").Value = "'49.27"
$ws.Range("E12").Value = "'  +1.05%  "
$ws.Range("D13").Value = "'0.0000286"
$ws.Range("E13").Value = "'  +2.93%  "
$ws.Range("D14").Value = "'4.303.69"
$ws.Range("E14").Value = "'  +8.59%  "
$ws.Range("D15").Value = "'685.01"
$ws.Range("E15").Value = "'  -0.90%  "
$ws.Range("E16").Value = "'  +4.64%  "
$ws.Range("D17").Value = "'3.708.60"
$ws.Range("E17").Value = "'  +8.44%  "
$ws.Range("D18").Value = "'71.662.57"
$ws.Range("E18").Value = "'  +3.39%  "
$ws.Range("E19").Value = "'  +1.42%  "
$ws.Range("D20").Value = "'17.99"
$ws.Range("E20").Value = "'  +2.01%  "
$ws.Range("D21").Value = "'11.60"
$ws.Range("E21").Value = "'  +2.21%  "
$ws.Range("E22").Value = "'  +18.51%  "
$ws.Range("D23").Value = "'0.943"
$ws.Range("E23").Value = "'  +5.44%  "
$ws.Range("D24").Value = "'17.48"
$ws.Range("E24").Value = "'  +3.50%  "
$ws.Range("D25").Value = "'102.62"
$ws.Range("E25").Value = "'  +2.06%  "
$ws.Range("D26").Value = "'3.99"
$ws.Range("E26").Value = "'  +3.04%  "
$ws.Range("D27").Value = "'2.83"
$ws.Range("E27").Value = "'  +6.64%  "
$ws.Range("D28").Value = "'10.38"
$ws.Range("E28").Value = "'  +8.74%  "
$ws.Range("D29").Value = "'35.33"
$ws.Range("E29").Value = "'  +6.04%  "
$ws.Range("D30").Value = "'9.18"
$ws.Range("E30").Value = "'  +5.31%  "
$ws.Range("E31").Value = "'  +5.31%  "
$ws.Range("E32").Value = "'  +11.94%  "
$ws.Range("D33").Value = "'592.50"
$ws.Range("E33").Value = "'  +4.14%  "
$ws.Range("D34").Value = "'11.22"
$ws.Range("E34").Value = "'  +2.30%  "
$ws.Range("E35").Value = "'  +4.76%  "
$ws.Range("D36").Value = "'59.13"
$ws.Range("E36").Value = "'  +1.76%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  -0.07%  "
$ws.Range("D38").Value = "'3.677.86"
$ws.Range("E38").Value = "'  +2.76%  "
$ws.Range("E39").Value = "'  +4.59%  "
$ws.Range("D40").Value = "'0.0₃0769"
$ws.Range("E40").Value = "'  +6.44%  "
$ws.Range("D41").Value = "'35.41"
$ws.Range("E41").Value = "'  +1.80%  "
$ws.Range("D42").Value = "'3.43"
$ws.Range("E42").Value = "'  +5.55%  "
$ws.Range("D43").Value = "'2.79"
$ws.Range("E43").Value = "'  +5.58%  "
$ws.Range("D44").Value = "'0.0458"
$ws.Range("E44").Value = "'  +10.19%  "
$ws.Range("E45").Value = "'  +5.41%  "
$ws.Range("B46").Value = "'ApeXProtocol"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.39"
$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("B47").Value = "'ThetaToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.88"
$ws.Range("E47").Value = "'  +9.34%  "
$ws.Range("E48").Value = "'  +4.15%  "
$ws.Range("E49").Value = "'  +0.45%  "
$ws.Range("E50").Value = "'  -0.16%  "
$ws.Range("D51").Value = "'135.99"
$ws.Range("E51").Value = "'  +3.08%  "

Write-Output "Updated crypto figures on $($ws.Name)"
